$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '27.696.59'
Set-TextValue $ws.Range("E2") '  +1.49%  '
Set-TextValue $ws.Range("D3") '1.868.10'
Set-TextValue $ws.Range("E3") '  +0.90%  '
Set-TextValue $ws.Range("D4") '1.008'
Set-TextValue $ws.Range("E4") '  +0.57%  '
Set-TextValue $ws.Range("D5") '333.19'
Set-TextValue $ws.Range("E5") '  +3.55%  '
Set-TextValue $ws.Range("D6") '1.008'
Set-TextValue $ws.Range("E6") '  +0.59%  '
Set-TextValue $ws.Range("D7") '0.4689'
Set-TextValue $ws.Range("E7") '  +4.27%  '
Set-TextValue $ws.Range("D8") '0.3938'
Set-TextValue $ws.Range("E8") '  +2.32%  '
Set-TextValue $ws.Range("D9") '47.38'
Set-TextValue $ws.Range("E9") '  -0.90%  '
Set-TextValue $ws.Range("D10") '0.08020'
Set-TextValue $ws.Range("E10") '  +2.08%  '
Set-TextValue $ws.Range("D11") '1.023'
Set-TextValue $ws.Range("E11") '  +1.01%  '
Set-TextValue $ws.Range("D12") '21.71'
Set-TextValue $ws.Range("E12") '  +2.16%  '
Set-TextValue $ws.Range("D13") '1.847.05'
Set-TextValue $ws.Range("E13") '  +0.42%  '
Set-TextValue $ws.Range("D14") '5.927'
Set-TextValue $ws.Range("E14") '  +1.19%  '
Set-TextValue $ws.Range("D15") '7.102'
Set-TextValue $ws.Range("E15") '  -0.19%  '
Set-TextValue $ws.Range("D16") '1.010'
Set-TextValue $ws.Range("E16") '  +0.65%  '
Set-TextValue $ws.Range("D17") '0.00001046'
Set-TextValue $ws.Range("E17") '  +1.72%  '
Set-TextValue $ws.Range("D18") '86.75'
Set-TextValue $ws.Range("E18") '  +1.74%  '
Set-TextValue $ws.Range("D19") '0.06649'
Set-TextValue $ws.Range("E19") '  +1.89%  '
Set-TextValue $ws.Range("D20") '17.24'
Set-TextValue $ws.Range("E20") '  +1.76%  '
Set-TextValue $ws.Range("D21") '1.008'
Set-TextValue $ws.Range("E21") '  +0.65%  '
Set-TextValue $ws.Range("D22") '27.717.92'
Set-TextValue $ws.Range("E22") '  +1.64%  '
Set-TextValue $ws.Range("D23") '5.466'
Set-TextValue $ws.Range("E23") '  -0.15%  '
Set-TextValue $ws.Range("D24") '10.99'
Set-TextValue $ws.Range("E24") '  +2.02%  '
Set-TextValue $ws.Range("D25") '2.314'
Set-TextValue $ws.Range("E25") '  +2.50%  '
Set-TextValue $ws.Range("D26") '2.098.00'
Set-TextValue $ws.Range("E26") '  +1.85%  '
Set-TextValue $ws.Range("D27") '157.90'
Set-TextValue $ws.Range("E27") '  +4.59%  '
Set-TextValue $ws.Range("D28") '20.12'
Set-TextValue $ws.Range("E28") '  +3.06%  '
Set-TextValue $ws.Range("D29") '2.092'
Set-TextValue $ws.Range("E29") '  +2.68%  '
Set-TextValue $ws.Range("D30") '5.559'
Set-TextValue $ws.Range("E30") '  +0.26%  '
Set-TextValue $ws.Range("D31") '122.16'
Set-TextValue $ws.Range("E31") '  +2.09%  '
Set-TextValue $ws.Range("D32") '0.9761'
Set-TextValue $ws.Range("E32") '  +4.55%  '
Set-TextValue $ws.Range("D33") '0.09525'
Set-TextValue $ws.Range("E33") '  +2.45%  '
Set-TextValue $ws.Range("D34") '1.451'
Set-TextValue $ws.Range("E34") '  -0.29%  '
Set-TextValue $ws.Range("D35") '3.612'
Set-TextValue $ws.Range("E35") '  +0.76%  '
Set-TextValue $ws.Range("D36") '5.299'
Set-TextValue $ws.Range("E36") '  +0.87%  '
Set-TextValue $ws.Range("D37") '0.02264'
Set-TextValue $ws.Range("E37") '  +1.95%  '
Set-TextValue $ws.Range("D38") '0.06079'
Set-TextValue $ws.Range("E38") '  +1.90%  '
Set-TextValue $ws.Range("D39") '1.232'
Set-TextValue $ws.Range("E39") '  +2.92%  '
Set-TextValue $ws.Range("D40") '8.148'
Set-TextValue $ws.Range("E40") '  -1.87%  '
Set-TextValue $ws.Range("B41") 'Frax'
Set-TextValue $ws.Range("C41") 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws.Range("D41") '1.007'
Set-TextValue $ws.Range("E41") '  +0.62%  '
Set-TextValue $ws.Range("B42") 'TheSandbox'
Set-TextValue $ws.Range("C42") 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range("D42") '0.5992'
Set-TextValue $ws.Range("E42") '  +1.69%  '
Set-TextValue $ws.Range("D43") '0.1900'
Set-TextValue $ws.Range("E43") '  +2.06%  '
Set-TextValue $ws.Range("D44") '10.26'
Set-TextValue $ws.Range("E44") '  +1.73%  '
Set-TextValue $ws.Range("D45") '1.261'
Set-TextValue $ws.Range("E45") '  +0.06%  '
Set-TextValue $ws.Range("D46") '0.5706'
Set-TextValue $ws.Range("E46") '  +0.86%  '
Set-TextValue $ws.Range("D47") '12.11'
Set-TextValue $ws.Range("E47") '  +2.08%  '
Set-TextValue $ws.Range("D48") '3.417'
Set-TextValue $ws.Range("E48") '  +1.66%  '
Set-TextValue $ws.Range("D49") '1.936'
Set-TextValue $ws.Range("E49") '  +0.33%  '
Set-TextValue $ws.Range("D50") '0.06831'
Set-TextValue $ws.Range("E50") '  -0.44%  '
Set-TextValue $ws.Range("D51") '113.23'
Set-TextValue $ws.Range("E51") '  +4.58%  '
